$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# --- Simple value replacements (rows whose position does not move) ---
$t.Cell(1,1).Range.Text = "0M"
$t.Cell(2,1).Range.Text = "0M"
$t.Cell(3,1).Range.Text = "0M"
$t.Cell(4,1).Range.Text = "1136"
# row 5 (0.00003) stays unchanged
$t.Cell(6,1).Range.Text = "0.00277"

# --- Remove the three rows that were dropped (0.00013 / 0.00005 / 0.00014) ---
# Delete from the highest index down so earlier indices stay valid.
$t.Rows.Item(9).Delete()
$t.Rows.Item(8).Delete()
$t.Rows.Item(7).Delete()

# After the deletions, row 10 (0.00018) shifted to 7 and stays unchanged;
# the following two rows need their values updated.
$t.Cell(8,1).Range.Text = "0.00007"
$t.Cell(9,1).Range.Text = "0.00028"

# --- Insert three new rows right after (former row 12, now row 9) ---
$refRow1 = $t.Rows.Item(10)
$t.Rows.Add($refRow1) | Out-Null
$t.Cell(10,1).Range.Text = "0.00040"

$refRow2 = $t.Rows.Item(11)
$t.Rows.Add($refRow2) | Out-Null
$t.Cell(11,1).Range.Text = "0.00049"

$refRow3 = $t.Rows.Item(12)
$t.Rows.Add($refRow3) | Out-Null
$t.Cell(12,1).Range.Text = "0.23914"

# --- Collapse the three multi-column (tab separated) summary rows down to a
#     single value each ---
$t.Cell(44,1).Range.Text = "99.82"
$t.Cell(45,1).Range.Text = "0.24"
$t.Cell(46,1).Range.Text = "134"
